# 0310 titrations and sample ID notation changes
# Adds a new titration data row (row 65) to the CRMAccuracyData sheet,
# extending the shared %-off formula down one row and scrolling the view.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CRMAccuracyData")

$newRow = 65

$ws.Cells.Item($newRow, 1).Value = 20220310
$ws.Cells.Item($newRow, 2).Value = 2221.9227812661102
$ws.Cells.Item($newRow, 3).Value = 2224.4699999999998
$ws.Cells.Item($newRow, 4).Formula = "=100*(B$newRow-C$newRow)/C$newRow"
$ws.Cells.Item($newRow, 5).Value = 180
$ws.Cells.Item($newRow, 6).Value = "CRM OPENED 20220302"

# Scroll/selection state matching the saved view
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 48
$excel.ActiveWindow.ScrollColumn = 1
[void]$ws.Range("D61").Select()
